# WEMUserDSN functions completed & tests updated
#
# This script updates the "Commands and Aliasses" worksheet:
#  - clears the AutoFilter criteria (column "State" filtered to "In Development")
#    and shows all previously filtered-out rows again
#  - updates the "State" value of the WEMUserDSN-related rows (A22 / A49)
#    from "In Development" to "Ready for testing"
#  - updates the active selection on the sheet to A49

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands and Aliasses")
$ws.Activate()

# The two WEMUserDSN command rows move from "In Development" to "Ready for testing"
$ws.Range("A22").Value = "Ready for testing"
$ws.Range("A49").Value = "Ready for testing"

# Clear the autofilter criteria on column "State" (removes the <filterColumn>
# that restricted the view to "In Development") and unhide every row that
# the filter had hidden.
$ws.ShowAllData()

# Leave the sheet with A49 selected, as it was after the edit.
$ws.Range("A49").Select()
